# Auto-generated Excel COM-interop edit script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets,
# reflecting refreshed market-board price data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value2 = 5267.2085
$ws.Range("J76").Value2 = 6130.615
$ws.Range("L76").Value2 = 6130.615
$ws.Range("N76").Value2 = -6760.615

# Row 79
$ws.Range("H79").Value2 = 5267.2085
$ws.Range("J79").Value2 = 6130.615
$ws.Range("L79").Value2 = 6130.615
$ws.Range("N79").Value2 = -8314.615

# Row 106
$ws.Range("H106").Value2 = 21959.75
$ws.Range("I106").Value2 = 21959.75
$ws.Range("K106").Value2 = 21959.75
$ws.Range("M106").Value2 = -21328.75

# Row 137
$ws.Range("H137").Value2 = 3727.875
$ws.Range("I137").Value2 = 2071.2856
$ws.Range("J137").Value2 = 5016.3335
$ws.Range("K137").Value2 = 6213.8568
$ws.Range("L137").Value2 = 15049.0005
$ws.Range("M137").Value2 = -3663.8568
$ws.Range("N137").Value2 = -20149.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value2 = 499.5
$ws.Range("I4").Value2 = 499.5
$ws.Range("K4").Value2 = 499.5
$ws.Range("M4").Value2 = -383.5

# Row 32
$ws.Range("H32").Value2 = 7238.675
$ws.Range("I32").Value2 = 5265.1943
$ws.Range("J32").Value2 = 25000
$ws.Range("K32").Value2 = 5265.1943
$ws.Range("L32").Value2 = 25000
$ws.Range("M32").Value2 = -4978.1943
$ws.Range("N32").Value2 = -25574

# Row 45
$ws.Range("H45").Value2 = 1534.3572
$ws.Range("I45").Value2 = 1506.8334
$ws.Range("J45").Value2 = 1699.5
$ws.Range("K45").Value2 = 1506.8334
$ws.Range("L45").Value2 = 1699.5
$ws.Range("M45").Value2 = -1129.8334
$ws.Range("N45").Value2 = -2453.5

# Row 74
$ws.Range("H74").Value2 = 3033.0667
$ws.Range("I74").Value2 = 1281.2
$ws.Range("J74").Value2 = 6536.8
$ws.Range("K74").Value2 = 1281.2
$ws.Range("L74").Value2 = 6536.8
$ws.Range("M74").Value2 = -407.2
$ws.Range("N74").Value2 = -8284.799999999999

# Row 77
$ws.Range("H77").Value2 = 3033.0667
$ws.Range("I77").Value2 = 1281.2
$ws.Range("J77").Value2 = 6536.8
$ws.Range("K77").Value2 = 6406
$ws.Range("L77").Value2 = 32684
$ws.Range("M77").Value2 = -2038
$ws.Range("N77").Value2 = -41420

# Row 97
$ws.Range("H97").Value2 = 987.6667
$ws.Range("I97").Value2 = 920.3333
$ws.Range("K97").Value2 = 920.3333
$ws.Range("M97").Value2 = -424.3333

# Row 102
$ws.Range("H102").Value2 = 1087.7333
$ws.Range("I102").Value2 = 863.125
$ws.Range("J102").Value2 = 1344.4286
$ws.Range("K102").Value2 = 863.125
$ws.Range("L102").Value2 = 1344.4286
$ws.Range("M102").Value2 = 758.875
$ws.Range("N102").Value2 = -4588.4286

# Row 122
$ws.Range("H122").Value2 = 419231.16
$ws.Range("I122").Value2 = 502427.9
$ws.Range("J122").Value2 = 3247.5
$ws.Range("K122").Value2 = 1507283.7
$ws.Range("L122").Value2 = 9742.5
$ws.Range("M122").Value2 = -1504833.7
$ws.Range("N122").Value2 = -14642.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value2 = 516
$ws.Range("I20").Value2 = 557.75
$ws.Range("K20").Value2 = 557.75
$ws.Range("M20").Value2 = -310.75

# Row 99
$ws.Range("H99").Value2 = 3154.913
$ws.Range("I99").Value2 = 2947.5293
$ws.Range("J99").Value2 = 3742.5
$ws.Range("K99").Value2 = 2947.5293
$ws.Range("L99").Value2 = 3742.5
$ws.Range("M99").Value2 = -1449.5293
$ws.Range("N99").Value2 = -6738.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 4198.6875
$ws.Range("I31").Value2 = 2772.6365
$ws.Range("J31").Value2 = 7336
$ws.Range("K31").Value2 = 2772.6365
$ws.Range("L31").Value2 = 7336
$ws.Range("M31").Value2 = -2477.6365
$ws.Range("N31").Value2 = -7926

# Row 34
$ws.Range("H34").Value2 = 4198.6875
$ws.Range("I34").Value2 = 2772.6365
$ws.Range("J34").Value2 = 7336
$ws.Range("K34").Value2 = 2772.6365
$ws.Range("L34").Value2 = 7336
$ws.Range("M34").Value2 = -2570.6365
$ws.Range("N34").Value2 = -7740

# Row 64
$ws.Range("H64").Value2 = 20666.334
$ws.Range("J64").Value2 = 20666.334
$ws.Range("L64").Value2 = 20666.334
$ws.Range("N64").Value2 = -21162.334

# Row 67
$ws.Range("H67").Value2 = 20666.334
$ws.Range("J67").Value2 = 20666.334
$ws.Range("L67").Value2 = 20666.334
$ws.Range("N67").Value2 = -22382.334

# Row 99
$ws.Range("H99").Value2 = 13917.333
$ws.Range("I99").Value2 = 9996.909
$ws.Range("K99").Value2 = 9996.909
$ws.Range("M99").Value2 = -8498.909

# Row 104
$ws.Range("H104").Value2 = 9750
$ws.Range("J104").Value2 = 9750
$ws.Range("L104").Value2 = 9750
$ws.Range("N104").Value2 = -14992

# Row 107
$ws.Range("H107").Value2 = 553.73914
$ws.Range("I107").Value2 = 375.05884
$ws.Range("J107").Value2 = 1060
$ws.Range("K107").Value2 = 375.05884
$ws.Range("L107").Value2 = 1060
$ws.Range("M107").Value2 = 1544.94116
$ws.Range("N107").Value2 = -4900

# Row 109
$ws.Range("H109").Value2 = 52753
$ws.Range("I109").Value2 = 40259
$ws.Range("J109").Value2 = 59000
$ws.Range("K109").Value2 = 40259
$ws.Range("L109").Value2 = 59000
$ws.Range("M109").Value2 = -39219
$ws.Range("N109").Value2 = -61080

# Row 126
$ws.Range("H126").Value2 = 13917.333
$ws.Range("I126").Value2 = 9996.909
$ws.Range("K126").Value2 = 29990.727
$ws.Range("M126").Value2 = -27520.727

# Row 132
$ws.Range("H132").Value2 = 1902.1333
$ws.Range("I132").Value2 = 1149.6666
$ws.Range("J132").Value2 = 4912
$ws.Range("K132").Value2 = 3448.9998
$ws.Range("L132").Value2 = 14736
$ws.Range("M132").Value2 = -918.9998000000001
$ws.Range("N132").Value2 = -19796

# Row 134
$ws.Range("H134").Value2 = 2691.913
$ws.Range("I134").Value2 = 2205.3157
$ws.Range("K134").Value2 = 6615.9471
$ws.Range("M134").Value2 = -4080.9471

# Row 138
$ws.Range("H138").Value2 = 117316.664
$ws.Range("J138").Value2 = 117316.664
$ws.Range("L138").Value2 = 117316.664
$ws.Range("N138").Value2 = -127596.664

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value2 = 194.5
$ws.Range("I60").Value2 = 151
$ws.Range("K60").Value2 = 453
$ws.Range("M60").Value2 = -202

# Row 105
$ws.Range("H105").Value2 = 0
$ws.Range("J105").Value2 = 0
$ws.Range("L105").Value2 = 0
$ws.Range("N105").ClearContents()

# Row 122
$ws.Range("H122").Value2 = 296.72726
$ws.Range("I122").Value2 = 268.3889
$ws.Range("J122").Value2 = 424.25
$ws.Range("K122").Value2 = 2415.5001
$ws.Range("L122").Value2 = 3818.25
$ws.Range("M122").Value2 = 34.49990000000025
$ws.Range("N122").Value2 = -8718.25

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value2 = 1902.6666
$ws.Range("I113").Value2 = 1902.6666
$ws.Range("K113").Value2 = 1902.6666
$ws.Range("M113").Value2 = 267.3334

# Row 122
$ws.Range("H122").Value2 = 61418.47
$ws.Range("I122").Value2 = 2274.3333
$ws.Range("K122").Value2 = 6822.999899999999
$ws.Range("M122").Value2 = -4372.999899999999

# Row 126
$ws.Range("H126").Value2 = 3796.4285
$ws.Range("I126").Value2 = 2679
$ws.Range("J126").Value2 = 4101.1816
$ws.Range("K126").Value2 = 8037
$ws.Range("L126").Value2 = 12303.5448
$ws.Range("M126").Value2 = -5567
$ws.Range("N126").Value2 = -17243.5448

# Row 132
$ws.Range("H132").Value2 = 1953.8
$ws.Range("I132").Value2 = 921.6667
$ws.Range("K132").Value2 = 2765.0001
$ws.Range("M132").Value2 = -235.0001000000002

# Row 134
$ws.Range("H134").Value2 = 95853.28999999999
$ws.Range("J134").Value2 = 95853.28999999999
$ws.Range("L134").Value2 = 287559.87
$ws.Range("N134").Value2 = -292629.87

# Row 141
$ws.Range("H141").Value2 = 78999.5
$ws.Range("J141").Value2 = 78999.5
$ws.Range("L141").Value2 = 78999.5
$ws.Range("N141").Value2 = -89359.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 1516
$ws.Range("I7").Value2 = 1754.875
$ws.Range("J7").Value2 = 1303.6666
$ws.Range("K7").Value2 = 1754.875
$ws.Range("L7").Value2 = 1303.6666
$ws.Range("M7").Value2 = -1642.875
$ws.Range("N7").Value2 = -1527.6666

# Row 40
$ws.Range("H40").Value2 = 1430.8125
$ws.Range("I40").Value2 = 1426.2
$ws.Range("K40").Value2 = 1426.2
$ws.Range("M40").Value2 = -1290.2

# Row 55
$ws.Range("H55").Value2 = 778.0833
$ws.Range("I55").Value2 = 736.8
$ws.Range("K55").Value2 = 736.8
$ws.Range("M55").Value2 = -563.8

# Row 100
$ws.Range("H100").Value2 = 4786.5
$ws.Range("J100").Value2 = 12850
$ws.Range("L100").Value2 = 12850
$ws.Range("N100").Value2 = -13932

# Row 126
$ws.Range("H126").Value2 = 1516
$ws.Range("I126").Value2 = 1754.875
$ws.Range("J126").Value2 = 1303.6666
$ws.Range("K126").Value2 = 5264.625
$ws.Range("L126").Value2 = 3910.9998
$ws.Range("M126").Value2 = -2794.625
$ws.Range("N126").Value2 = -8850.9998

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value2 = 18500
$ws.Range("J15").Value2 = 18500
$ws.Range("L15").Value2 = 18500
$ws.Range("N15").Value2 = -19076

# Row 107
$ws.Range("H107").Value2 = 603.4666999999999
$ws.Range("I107").Value2 = 393.125
$ws.Range("J107").Value2 = 843.8570999999999
$ws.Range("K107").Value2 = 1179.375
$ws.Range("L107").Value2 = 2531.5713
$ws.Range("M107").Value2 = 740.625
$ws.Range("N107").Value2 = -6371.5713

# Row 122
$ws.Range("H122").Value2 = 6934.533
$ws.Range("I122").Value2 = 6216.2856
$ws.Range("K122").Value2 = 18648.8568
$ws.Range("M122").Value2 = -16198.8568

# Row 126
$ws.Range("H126").Value2 = 1944.6875
$ws.Range("I126").Value2 = 1640.6666
$ws.Range("K126").Value2 = 4921.9998
$ws.Range("M126").Value2 = -2451.9998

# Row 132
$ws.Range("H132").Value2 = 988
$ws.Range("I132").Value2 = 988
$ws.Range("K132").Value2 = 2964
$ws.Range("M132").Value2 = -434

# Row 140
$ws.Range("H140").Value2 = 96666.336
$ws.Range("J140").Value2 = 124999.5
$ws.Range("L140").Value2 = 124999.5
$ws.Range("N140").Value2 = -135359.5
